$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force-text cells: values that look numeric must stay as text strings,
# matching the original inlineStr storage. Apply a Text number format
# before writing so Excel does not auto-convert them to numbers.

$ws.Range("D2").Value = "63.488.26"
$ws.Range("E2").Value = "  +4.97%  "

$ws.Range("D3").Value = "3.059.61"
$ws.Range("E3").Value = "  +3.20%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "548.97"
$ws.Range("E5").Value = "  +5.04%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.11"
$ws.Range("E6").Value = "  +7.88%  "

$ws.Range("E7").Value = "  +0.12%  "

$ws.Range("D8").Value = "3.052.80"
$ws.Range("E8").Value = "  +3.17%  "

$ws.Range("E9").Value = "  +3.30%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.24"
$ws.Range("E10").Value = "  +2.88%  "

$ws.Range("E11").Value = "  +1.87%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.452"
$ws.Range("E12").Value = "  +4.25%  "

$ws.Range("E13").Value = "  +5.10%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.63"
$ws.Range("E14").Value = "  +5.31%  "

$ws.Range("D15").Value = "3.567.01"
$ws.Range("E15").Value = "  +3.77%  "

$ws.Range("D16").Value = "63.528.46"
$ws.Range("E16").Value = "  +5.09%  "

$ws.Range("D17").Value = "3.061.29"
$ws.Range("E17").Value = "  +3.38%  "

$ws.Range("E18").Value = "  -0.75%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.71"
$ws.Range("E19").Value = "  +4.35%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "477.42"
$ws.Range("E20").Value = "  +5.50%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.54"
$ws.Range("E21").Value = "  +4.81%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.677"
$ws.Range("E22").Value = "  +2.35%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.19"
$ws.Range("E23").Value = "  +6.88%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.15"
$ws.Range("E24").Value = "  +4.61%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.51"
$ws.Range("E25").Value = "  +7.69%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.08%  "

$ws.Range("E27").Value = "  +5.49%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.94"
$ws.Range("E28").Value = "  +5.36%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.98"
$ws.Range("E29").Value = "  +9.11%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.26%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.96"
$ws.Range("E31").Value = "  +4.44%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.14"
$ws.Range("E32").Value = "  +1.06%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.41"
$ws.Range("E33").Value = "  +8.77%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.65"
$ws.Range("E34").Value = "  +8.07%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "55.79"
$ws.Range("E35").Value = "  +2.19%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.97"
$ws.Range("E36").Value = "  +4.72%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "465.08"
$ws.Range("E37").Value = "  +4.00%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0813"
$ws.Range("E38").Value = "  +6.11%  "

$ws.Range("D39").Value = "3.127.82"
$ws.Range("E39").Value = "  -0.76%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0394"
$ws.Range("E40").Value = "  +5.51%  "

$ws.Range("E41").Value = "  +3.54%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.21"
$ws.Range("E42").Value = "  +3.68%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.56"
$ws.Range("E43").Value = "  +8.31%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "27.93"
$ws.Range("E44").Value = "  +12.33%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.251"
$ws.Range("E45").Value = "  +4.83%  "

$ws.Range("E46").Value = "  -0.10%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.05"
$ws.Range("E47").Value = "  +7.39%  "

$ws.Range("E48").Value = "  +2.35%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₃0509"
$ws.Range("E49").Value = "  +3.78%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "115.73"
$ws.Range("E50").Value = "  -1.02%  "

$ws.Range("E51").Value = "  +7.11%  "

